$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Totale Positivi" (H) and "Tamponi" (I) counts for rows 2-6,
# matching the notebook change that dropped the swab counts for those days.
$ws.Range("H2:I6").ClearContents()
